# Rename "Sheet1" to "S3_Overview" and make it the active/selected sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "S3_Overview"

$ws.Select()
$ws.Activate()
